$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.856696666666667
$ws.Range("H2").Value = 8.57009
$ws.Range("I2").Value = 0.05747862151401942
$ws.Range("J2").Value = 0.05747862151401942
$ws.Range("O2").Value = 0.7792485920506572
$ws.Range("P2").Value = 0.7792485920506572
$ws.Range("Q2").Value = 0.1281990240777778
$ws.Range("R2").Value = 1.1537912167
$ws.Range("S2").Value = 0.04479013488781225
$ws.Range("T2").Value = 0.04479013488781225
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.856696666666667
$ws.Range("H3").Value = 8.57009
$ws.Range("I3").Value = 0.05747862151401942
$ws.Range("J3").Value = 0.05747862151401942
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.012713
$ws.Range("N3").Value = 0.03813900000000001
$ws.Range("O3").Value = 0.2207514079493428
$ws.Range("P3").Value = 0.2207514079493428
$ws.Range("Q3").Value = 0.03631718472333334
$ws.Range("R3").Value = 0.3268546625100001
$ws.Range("S3").Value = 0.01268848662620717
$ws.Range("T3").Value = 0.01268848662620717
$ws.Range("I4").Value = 0.7708435061432634
$ws.Range("J4").Value = 0.7708435061432632
$ws.Range("O4").Value = 0.7792485920506572
$ws.Range("P4").Value = 0.7792485920506572
$ws.Range("S4").Value = 0.6006787168535301
$ws.Range("T4").Value = 0.60067871685353
$ws.Range("I5").Value = 0.7708435061432634
$ws.Range("J5").Value = 0.7708435061432632
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.012713
$ws.Range("N5").Value = 0.03813900000000001
$ws.Range("O5").Value = 0.2207514079493428
$ws.Range("P5").Value = 0.2207514079493428
$ws.Range("Q5").Value = 0.4870483193226667
$ws.Range("R5").Value = 4.383434873904001
$ws.Range("S5").Value = 0.1701647892897332
$ws.Range("T5").Value = 0.1701647892897332
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 2.081608
$ws.Range("H6").Value = 6.244823999999999
$ws.Range("I6").Value = 0.04188332620983732
$ws.Range("J6").Value = 0.04188332620983732
$ws.Range("O6").Value = 0.7792485920506572
$ws.Range("P6").Value = 0.7792485920506572
$ws.Range("Q6").Value = 0.09341562834666665
$ws.Range("R6").Value = 0.84074065512
$ws.Range("S6").Value = 0.03263752297941412
$ws.Range("T6").Value = 0.03263752297941412
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 2.081608
$ws.Range("H7").Value = 6.244823999999999
$ws.Range("I7").Value = 0.04188332620983732
$ws.Range("J7").Value = 0.04188332620983732
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.012713
$ws.Range("N7").Value = 0.03813900000000001
$ws.Range("O7").Value = 0.2207514079493428
$ws.Range("P7").Value = 0.2207514079493428
$ws.Range("Q7").Value = 0.026463482504
$ws.Range("R7").Value = 0.238171342536
$ws.Range("S7").Value = 0.009245803230423198
$ws.Range("T7").Value = 0.009245803230423198
$ws.Range("G8").Value = 4.573220666666667
$ws.Range("H8").Value = 13.719662
$ws.Range("I8").Value = 0.09201621679565497
$ws.Range("J8").Value = 0.09201621679565496
$ws.Range("O8").Value = 0.7792485920506572
$ws.Range("P8").Value = 0.7792485920506572
$ws.Range("Q8").Value = 0.2052308994511111
$ws.Range("R8").Value = 1.84707809506
$ws.Range("S8").Value = 0.07170350738384217
$ws.Range("T8").Value = 0.07170350738384217
$ws.Range("G9").Value = 4.573220666666667
$ws.Range("H9").Value = 13.719662
$ws.Range("I9").Value = 0.09201621679565497
$ws.Range("J9").Value = 0.09201621679565496
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.012713
$ws.Range("N9").Value = 0.03813900000000001
$ws.Range("O9").Value = 0.2207514079493428
$ws.Range("P9").Value = 0.2207514079493428
$ws.Range("Q9").Value = 0.05813935433533334
$ws.Range("R9").Value = 0.523254189018
$ws.Range("S9").Value = 0.0203127094118128
$ws.Range("T9").Value = 0.02031270941181279
$ws.Range("G10").Value = 1.877589
$ws.Range("H10").Value = 5.632767
$ws.Range("I10").Value = 0.037778329337225
$ws.Range("J10").Value = 0.037778329337225
$ws.Range("O10").Value = 0.7792485920506572
$ws.Range("P10").Value = 0.7792485920506572
$ws.Range("Q10").Value = 0.08425993569000001
$ws.Range("R10").Value = 0.7583394212100001
$ws.Range("S10").Value = 0.02943870994605862
$ws.Range("T10").Value = 0.02943870994605862
$ws.Range("G11").Value = 1.877589
$ws.Range("H11").Value = 5.632767
$ws.Range("I11").Value = 0.037778329337225
$ws.Range("J11").Value = 0.037778329337225
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.012713
$ws.Range("N11").Value = 0.03813900000000001
$ws.Range("O11").Value = 0.2207514079493428
$ws.Range("P11").Value = 0.2207514079493428
$ws.Range("Q11").Value = 0.02386978895700001
$ws.Range("R11").Value = 0.214828100613
$ws.Range("S11").Value = 0.008339619391166382
$ws.Range("T11").Value = 0.00833961939116638
